$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.316.71'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '1.665.26'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5352'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2659'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06421'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07849'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.565'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.58%  '
$ws.Range("D13").Value = '1.658.62'
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '1.893.43'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5534'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").Value = '0.0₅8239'
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("E18").Value = '  +0.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.704'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.039'
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '146.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1232'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.198'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.484'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05837'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.621'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.280'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.619'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9699'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.825'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.422'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5825'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01606'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8718'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.868'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '105.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("D42").Value = '1.052.58'
$ws.Range("E42").Value = '  +2.72%  '
$ws.Range("D44").Value = '1.804.75'
$ws.Range("E44").Value = '  +0.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.87'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("D46").Value = '0.0₈106'
$ws.Range("E46").Value = '  -4.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.014'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.050'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.31%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.415'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.39%  '
